# Update the IA-vs-IA results grid: a handful of cells move from one
# "bucket" (0,1,2,3) to another. Each bucket has its own fill colour,
# driven directly by a cell style (no conditional-formatting rules are
# defined on this sheet), so every value change must be paired with the
# matching style change:
#   value 0 -> style with fill FF6666 (red)
#   value 1 -> style with fill FFC966 (orange)
#   value 2 -> style with fill 6DC066 (green)
#   value 3 -> style with fill 8067A2 (purple)
#
# Setting .Interior.Color directly would synthesize a brand-new style/fill
# entry in the workbook instead of reusing the existing one, so instead we
# copy the *format only* from a stable donor cell that already carries the
# right style, then overwrite the value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stable donor cells (untouched by this edit) for each of the four buckets.
$donorForStyle = @{
    2 = "B4"   # value 1 -> orange fill
    3 = "L4"   # value 0 -> red fill
    4 = "D5"   # value 2 -> green fill
    5 = "I33"  # value 3 -> purple fill
}

$changes = @(
    @{ Cell = "J5";  Style = 2; Value = 1 },
    @{ Cell = "C6";  Style = 4; Value = 2 },
    @{ Cell = "J7";  Style = 2; Value = 1 },
    @{ Cell = "B8";  Style = 3; Value = 0 },
    @{ Cell = "C8";  Style = 4; Value = 2 },
    @{ Cell = "I8";  Style = 2; Value = 1 },
    @{ Cell = "C9";  Style = 4; Value = 2 },
    @{ Cell = "E10"; Style = 4; Value = 2 },
    @{ Cell = "D11"; Style = 2; Value = 1 },
    @{ Cell = "K12"; Style = 2; Value = 1 },
    @{ Cell = "K13"; Style = 2; Value = 1 },
    @{ Cell = "G19"; Style = 2; Value = 1 },
    @{ Cell = "F21"; Style = 2; Value = 1 },
    @{ Cell = "G21"; Style = 3; Value = 0 },
    @{ Cell = "B22"; Style = 2; Value = 1 },
    @{ Cell = "F22"; Style = 2; Value = 1 },
    @{ Cell = "E23"; Style = 4; Value = 2 },
    @{ Cell = "F23"; Style = 4; Value = 2 },
    @{ Cell = "F24"; Style = 2; Value = 1 },
    @{ Cell = "G27"; Style = 2; Value = 1 },
    @{ Cell = "H33"; Style = 2; Value = 1 },
    @{ Cell = "B34"; Style = 2; Value = 1 },
    @{ Cell = "F34"; Style = 4; Value = 2 },
    @{ Cell = "H34"; Style = 2; Value = 1 },
    @{ Cell = "B35"; Style = 4; Value = 2 },
    @{ Cell = "B36"; Style = 2; Value = 1 },
    @{ Cell = "D36"; Style = 2; Value = 1 },
    @{ Cell = "G36"; Style = 2; Value = 1 },
    @{ Cell = "J36"; Style = 3; Value = 0 },
    @{ Cell = "D37"; Style = 5; Value = 3 },
    @{ Cell = "F37"; Style = 4; Value = 2 },
    @{ Cell = "B38"; Style = 2; Value = 1 },
    @{ Cell = "C38"; Style = 5; Value = 3 },
    @{ Cell = "D38"; Style = 2; Value = 1 },
    @{ Cell = "F38"; Style = 4; Value = 2 },
    @{ Cell = "G38"; Style = 2; Value = 1 },
    @{ Cell = "D39"; Style = 2; Value = 1 },
    @{ Cell = "F40"; Style = 2; Value = 1 },
    @{ Cell = "I40"; Style = 5; Value = 3 }
)

foreach ($change in $changes) {
    $donor = $donorForStyle[$change.Style]
    $ws.Range($donor).Copy()
    $ws.Range($change.Cell).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws.Range($change.Cell).Value = $change.Value
}

$excel.CutCopyMode = $false
